$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their text formatting (values include
# multi-part numbers such as "33.699.94" that Excel would otherwise reinterpret).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values from the refreshed crypto price feed
$ws.Range("D2").Value = '33.699.94'
$ws.Range("E2").Value = '  +6.74%  '

$ws.Range("D3").Value = '1.775.95'
$ws.Range("E3").Value = '  +3.89%  '

$ws.Range("D5").Value = '224.27'
$ws.Range("E5").Value = '  +0.68%  '

$ws.Range("D6").Value = '0.557'
$ws.Range("E6").Value = '  +3.74%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '30.15'
$ws.Range("E8").Value = '  +0.86%  '

$ws.Range("D9").Value = '46.52'
$ws.Range("E9").Value = '  +3.41%  '

$ws.Range("E10").Value = '  +2.97%  '

$ws.Range("E11").Value = '  +1.46%  '

$ws.Range("E12").Value = '  +1.29%  '

$ws.Range("D13").Value = '2.030.55'
$ws.Range("E13").Value = '  +3.85%  '

$ws.Range("D14").Value = '1.772.64'
$ws.Range("E14").Value = '  +3.69%  '

$ws.Range("E15").Value = '  +1.23%  '

$ws.Range("D16").Value = '33.668.91'
$ws.Range("E16").Value = '  +6.73%  '

$ws.Range("D17").Value = '10.04'
$ws.Range("E17").Value = '  -1.06%  '

$ws.Range("E18").Value = '  -0.69%  '

$ws.Range("D19").Value = '68.31'
$ws.Range("E19").Value = '  +1.39%  '

$ws.Range("D20").Value = '249.80'
$ws.Range("E20").Value = '  -0.47%  '

$ws.Range("E21").Value = '  +1.54%  '

$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").Value = '10.26'
$ws.Range("E23").Value = '  +1.20%  '

$ws.Range("E24").Value = '  -2.52%  '

$ws.Range("E25").Value = '  -2.33%  '

$ws.Range("D26").Value = '158.47'
$ws.Range("E26").Value = '  -0.60%  '

$ws.Range("D27").Value = '16.38'
$ws.Range("E27").Value = '  +1.82%  '

$ws.Range("E28").Value = '  +0.79%  '

$ws.Range("E29").Value = '  +2.05%  '

$ws.Range("E30").Value = '  -0.06%  '

$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("E32").Value = '  +2.18%  '

$ws.Range("E33").Value = '  +3.06%  '

$ws.Range("E34").Value = '  +3.92%  '

$ws.Range("E35").Value = '  +5.00%  '

$ws.Range("D36").Value = '1.481.37'
$ws.Range("E36").Value = '  -3.34%  '

$ws.Range("E37").Value = '  +2.12%  '

$ws.Range("D38").Value = '0.627'
$ws.Range("E38").Value = '  +2.12%  '

$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").Value = '83.06'
$ws.Range("E39").Value = '  +0.43%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.0184'
$ws.Range("E40").Value = '  +1.81%  '

$ws.Range("E41").Value = '  +1.41%  '

$ws.Range("D42").Value = '2.70'
$ws.Range("E42").Value = '  -0.37%  '

$ws.Range("E43").Value = '  +3.85%  '

$ws.Range("E44").Value = '  +0.86%  '

$ws.Range("D45").Value = '0.0509'
$ws.Range("E45").Value = '  +0.96%  '

$ws.Range("E46").Value = '  +4.33%  '

$ws.Range("D47").Value = '1.922.58'
$ws.Range("E47").Value = '  +3.93%  '

$ws.Range("E48").Value = '  -0.01%  '

$ws.Range("E49").Value = '  +1.46%  '

$ws.Range("E50").Value = '  +12.38%  '

$ws.Range("D51").Value = '50.86'
$ws.Range("E51").Value = '  -2.96%  '
